$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on columns B, C, D, E for the rows we touch so that
# values like dates/numbers are not auto-converted by Excel, matching the
# original inline-string (text) cell contents.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.908.79'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.816.11'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.91%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.79'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.03%  '

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.09%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4669'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.71%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.13%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07359'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.11%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8709'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.39%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.39'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.02%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.786.03'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +6.64%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.376'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.94%  '

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07069'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.27%  '

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.513'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.33%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.57'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.56%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.17%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008695'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.71%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.05%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.72'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.45%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.927.16'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.52%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.321'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.52%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.61'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.54%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.025.44'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +5.71%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.894'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.50%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '150.26'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.05%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.172'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.25%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.67%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.327'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.72%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.84'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.42%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08935'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.37%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7685'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.91%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.164'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.03%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.503'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.47%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.903'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.43%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.001'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.18%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.01%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01962'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.85%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05283'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.30%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.930'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.89%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.248'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.94%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5324'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.07%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.350'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.61%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.38%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.439'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.88%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4925'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.89%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.48'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.62%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.001'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.14%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.22%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '103.82'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.39%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06288'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.02%  '
